# Shift the weekly forecast dates forward by one week on the
# "Forecast Comparison" sheet (B2:B17), and correspondingly update the
# derived date references on the "Summary" sheet (Historical Range end
# date, Max/Min Forecast Week).
#
# The Week_Start_Date / *_Week cells are stored as plain text (not real
# Excel dates) in the source workbook, so for each cell we temporarily
# force a Text number format before assigning the string, then clear the
# format again so no stray formatting is left behind (Excel would
# otherwise auto-detect the "yyyy-mm-dd" string as a date literal and
# convert it to a date serial number).

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$newDates = @{
    2  = "2025-01-12"
    3  = "2025-01-19"
    4  = "2025-01-26"
    5  = "2025-02-02"
    6  = "2025-02-09"
    7  = "2025-02-16"
    8  = "2025-02-23"
    9  = "2025-03-02"
    10 = "2025-03-09"
    11 = "2025-03-16"
    12 = "2025-03-23"
    13 = "2025-03-30"
    14 = "2025-04-06"
    15 = "2025-04-13"
    16 = "2025-04-20"
    17 = "2025-04-27"
}

foreach ($row in $newDates.Keys) {
    $cell = $wsForecast.Range("B$row")
    Set-TextValue $cell $newDates[$row]
}

$wsSummary = $wb.Worksheets.Item("Summary")

Set-TextValue $wsSummary.Range("B2") "2023-01-01 to 2025-01-05"
Set-TextValue $wsSummary.Range("B13") "2025-01-12"
Set-TextValue $wsSummary.Range("B15") "2025-01-12"
